$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below row 6 (i.e. at row 7), shifting existing rows 7.. down.
$ws.Rows.Item(7).Insert()

# Copy the time value + formatting pattern from row 6's A cell into the new row 7's A cell.
$ws.Cells.Item(7, 1).Value = $ws.Cells.Item(6, 1).Value
$ws.Cells.Item(7, 1).NumberFormat = $ws.Cells.Item(6, 1).NumberFormat
$ws.Cells.Item(7, 1).Font.Bold = $ws.Cells.Item(6, 1).Font.Bold
$ws.Cells.Item(7, 1).HorizontalAlignment = $ws.Cells.Item(6, 1).HorizontalAlignment
$ws.Cells.Item(7, 1).VerticalAlignment = $ws.Cells.Item(6, 1).VerticalAlignment

$ws.Cells.Item(7, 2).Value = "midMarch 2025"

$ws.Range("B8").Select()
